$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.519.39'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.992.19'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +5.97%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4679'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.73%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3946'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.65%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.37'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07944'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.12%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.002'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.51%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.93'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.20%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.017.73'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +10.16%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.275'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.88%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.864'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.83%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.07141'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.52%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.73'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.002'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.06%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000009943'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.38'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.23%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '29.640.38'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.52%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.526'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.68%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.102'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.76'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.968'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '120.31'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.959'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09451'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9025'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.350'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.75%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.253'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('E35').Value = '  -2.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05831'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.176'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02119'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.30%  '
$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.000003320'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +102.51%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.866'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.91%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5752'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1825'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.797'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.05'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5369'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.41%  '
$ws.Range('E46').Value = '  +6.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.176'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.26%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06943'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.44%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.867'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '114.01'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3099'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.52%  '
